$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data range (rows 1-18, cols A:E) before rewriting content
$ws.Range("A1:E18").ClearContents() | Out-Null

# Write header + data rows with new import model content
$ws.Cells.Item(1,1).Value = "Comunidad autónoma"
$ws.Cells.Item(1,2).Value = "Provincia"
$ws.Cells.Item(2,1).Value = "Cataluña"
$ws.Cells.Item(2,2).Value = "Barcelona"
$ws.Cells.Item(3,1).Value = "Cataluña"
$ws.Cells.Item(3,2).Value = "Gerona"
$ws.Cells.Item(4,1).Value = "Cataluña"
$ws.Cells.Item(4,2).Value = "Lérida"
$ws.Cells.Item(5,1).Value = "Cataluña"
$ws.Cells.Item(5,2).Value = "Tarragona"
$ws.Cells.Item(6,1).Value = "Andalucía"
$ws.Cells.Item(6,2).Value = "Almeria"
$ws.Cells.Item(7,1).Value = "Andalucía"
$ws.Cells.Item(7,2).Value = "Cadiz"
$ws.Cells.Item(8,1).Value = "Andalucía"
$ws.Cells.Item(8,2).Value = "Sevilla"
$ws.Cells.Item(9,1).Value = "Andalucía"
$ws.Cells.Item(9,2).Value = "Cordoba"
$ws.Cells.Item(10,1).Value = "Canarias"
$ws.Cells.Item(10,2).Value = "Tenerife"
$ws.Cells.Item(11,1).Value = "Canarias"
$ws.Cells.Item(11,2).Value = "Fuerteventura"
$ws.Cells.Item(12,1).Value = "Canarias"
$ws.Cells.Item(12,2).Value = "Lanzarote"
$ws.Cells.Item(13,1).Value = "Canarias"
$ws.Cells.Item(13,2).Value = "Gran Canaria"
$ws.Cells.Item(14,1).Value = "Aragón"
$ws.Cells.Item(14,2).Value = "Huesca"
$ws.Cells.Item(15,1).Value = "Aragón"
$ws.Cells.Item(15,2).Value = "Teruel"
$ws.Cells.Item(16,1).Value = "Aragón"
$ws.Cells.Item(16,2).Value = "Zaragoza"
$ws.Cells.Item(17,1).Value = "Galicia"
$ws.Cells.Item(17,2).Value = "La Coruña"
$ws.Cells.Item(18,1).Value = "Galicia"
$ws.Cells.Item(18,2).Value = "Lugo"
$ws.Cells.Item(19,1).Value = "Galicia"
$ws.Cells.Item(19,2).Value = "Orense"
$ws.Cells.Item(20,1).Value = "Galicia"
$ws.Cells.Item(20,2).Value = "Pontevedra"

# Column A width (best-fit for longest label "Comunidad autónoma")
$ws.Columns.Item(1).ColumnWidth = 19

# Selection moved to E16 as left by the editing session
$ws.Range("E16").Select() | Out-Null
